# Add two new columns (I = "I0", J = "IF") to the stats sheet, matching
# the existing header style (bold, centered, thin border) used by the
# other header cells (e.g. H1), and fill in the per-row numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from an existing header cell onto the two
# new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Row -> (I value, J value)
$values = @{
    2  = @(7, 7)
    3  = @(8, 9)
    4  = @(6, 7)
    5  = @(5, 6)
    6  = @(7, 8)
    7  = @(8, 8)
    8  = @(7, 7)
    9  = @(10, 10)
    10 = @(6, 8)
    11 = @(1, 2)
    12 = @(7, 8)
    13 = @(6, 6)
    14 = @(9, 9)
    15 = @(8, 8)
    16 = @(5, 7)
    17 = @(6, 6)
    18 = @(5, 6)
    19 = @(7, 7)
    20 = @(9, 9)
    21 = @(4, 5)
    22 = @(9, 9)
    23 = @(6, 6)
    24 = @(6, 6)
    25 = @(3, 3)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}

Write-Output "Added I0/IF columns"
